# Generate Report for Archive
#
# 1) Every "Status" cell that currently reads "Ready for handoff" becomes
#    "In Translation" (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2) The "Status" column is narrowed on all three sheets (was ~17.216 chars
#    of OOXML column width, now ~13.410).

$wb = $excel.ActiveWorkbook

# --- 1) Update the status text wherever it appears -------------------------

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- 2) Narrow the status columns ------------------------------------------
# The underlying engine snaps ColumnWidth to its own pixel grid, so this is
# the closest settable width to the target OOXML width (13.4101845877511).
$newWidth = 12.4583333333333

$overview.Columns.Item(5).ColumnWidth = $newWidth   # column E
$overview.Columns.Item(6).ColumnWidth = $newWidth   # column F

$zhcn.Columns.Item(3).ColumnWidth = $newWidth        # column C

$dede.Columns.Item(3).ColumnWidth = $newWidth        # column C
